# Generate Report for Handoff
# Updates the "Latest Handoff Datetime" values for the
# fd237789-b6ac-4be7-a78e-86c2c38c8006 row (last row, row 7) across the
# zh-cn and de-de detail sheets, and the rolled-up "Latest Handoff Date"
# on the Overview sheet (which shares the de-de timestamp text).

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")
$overview = $wb.Worksheets.Item("Overview")

# zh-cn: Latest Handoff Datetime for fd237789-... moves from 23:02:48 to 23:03:05
$zhcn.Range("E7").Value = "2016-03-13 23:03:05"

# de-de: Latest Handoff Datetime for fd237789-... moves from 23:02:52 to 23:03:09
$dede.Range("E7").Value = "2016-03-13 23:03:09"

# Overview: Latest Handoff Date for fd237789-... mirrors the de-de timestamp
$overview.Range("D7").Value = "2016-03-13 23:03:09"
